$wb = $excel.ActiveWorkbook

# For each sheet, the most recent date's block of group rows (at the bottom of the
# sheet) gets each row quadruplicated (4 consecutive copies of each original row),
# pushing the block's total row count to 4x its original size. We process rows
# bottom-to-top within each block so inserts never disturb rows not yet handled.

function Expand-Block {
    param(
        [string]$SheetName,
        [int]$StartRow,
        [int]$EndRow,
        [int]$LastCol
    )

    $ws = $wb.Worksheets.Item($SheetName)
    $colLetter = [char](64 + $LastCol)

    for ($r = $EndRow; $r -ge $StartRow; $r--) {
        # Make room: insert 3 blank rows right below row $r
        $ws.Rows.Item($r + 1).Insert()
        $ws.Rows.Item($r + 1).Insert()
        $ws.Rows.Item($r + 1).Insert()

        # Snapshot the original row's values+formatting, then stamp 3 copies below it
        $srcRange = $ws.Range("A" + $r + ":" + $colLetter + $r)
        $vals = $srcRange.Value()

        $ws.Range("A" + ($r + 1) + ":" + $colLetter + ($r + 1)).Value = $vals
        $ws.Range("A" + ($r + 2) + ":" + $colLetter + ($r + 2)).Value = $vals
        $ws.Range("A" + ($r + 3) + ":" + $colLetter + ($r + 3)).Value = $vals
    }
}

# age: last block (date 2020-07-06) = rows 448..461 (14 group rows), columns A:F
Expand-Block "age" 448 461 6

# gender: last block (date 2020-07-06) = rows 108..111 (4 group rows), columns A:F
Expand-Block "gender" 108 111 6

# race: last block (date 2020-07-06) = rows 210..216 (7 group rows), columns A:F
Expand-Block "race" 210 216 6
